$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 50
$ws.Range("A50").Value = 44523
$ws.Range("A50").NumberFormat = 'm"月"d"日"'
$ws.Range("B50").Value = 1712
$ws.Range("C50").Value = "long"
$ws.Range("D50").Value = 215
$ws.Range("E50").Value = 27.75

# Row 51
$ws.Range("A51").Value = 44523
$ws.Range("A51").NumberFormat = 'm"月"d"日"'
$ws.Range("B51").Value = 2436
$ws.Range("C51").Value = "short"
$ws.Range("D51").Value = -60
$ws.Range("E51").Value = 101

# Update sheet view: scroll so row 46 is at the top and select F51 (matches the
# author's final cursor position after appending the new trades)
$win = $excel.ActiveWindow
$win.ScrollColumn = 1
$win.ScrollRow = 46
$null = $ws.Range("F51").Select()
